$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Replace the generic "Member 1".."Member 5" placeholders with the
# actual team member names.
$ws.Range("A8").Value = "Kunaal Sikka"
$ws.Range("A9").Value = "Mina Huh"
$ws.Range("A10").Value = "Vu Nguyen"
$ws.Range("A11").Value = "Nicolas Carmody"
$ws.Range("A12").Value = "Jonas Bokstaller"

# Update the current selection to match the saved state of the sheet.
$ws.Range("A8").Select()
